$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2116.3333
$ws.Range("I28").Value = 2140.4
$ws.Range("K28").Value = 2140.4
$ws.Range("M28").Value = -1655.4

$ws.Range("H38").Value = 1308.8182
$ws.Range("I38").Value = 112.125
$ws.Range("K38").Value = 336.375
$ws.Range("M38").Value = 35.625

$ws.Range("H41").Value = 2865.55
$ws.Range("I41").Value = 2982.158
$ws.Range("J41").Value = 650
$ws.Range("K41").Value = 2982.158
$ws.Range("L41").Value = 650
$ws.Range("M41").Value = -2542.158
$ws.Range("N41").Value = -1530

$ws.Range("H86").Value = 4531.75
$ws.Range("J86").Value = 5042.909
$ws.Range("L86").Value = 5042.909
$ws.Range("N86").Value = -7288.909

$ws.Range("H89").Value = 4531.75
$ws.Range("J89").Value = 5042.909
$ws.Range("L89").Value = 25214.545
$ws.Range("N89").Value = -36446.545

$ws.Range("H92").Value = 7210.467
$ws.Range("I92").Value = 8749.833000000001
$ws.Range("K92").Value = 8749.833000000001
$ws.Range("M92").Value = -7501.833000000001

$ws.Range("H98").Value = 731.82355
$ws.Range("I98").Value = 731.82355
$ws.Range("K98").Value = 731.82355
$ws.Range("M98").Value = 766.17645

$ws.Range("H99").Value = 1988.68
$ws.Range("I99").Value = 495
$ws.Range("J99").Value = 3162.2856
$ws.Range("K99").Value = 1485
$ws.Range("L99").Value = 9486.856800000001
$ws.Range("M99").Value = 13
$ws.Range("N99").Value = -12482.8568

$ws.Range("H122").Value = 731.82355
$ws.Range("I122").Value = 731.82355
$ws.Range("K122").Value = 2195.47065
$ws.Range("M122").Value = 254.5293500000002

$ws.Range("H131").Value = 6457.8335
$ws.Range("I131").Value = 2499.1667
$ws.Range("J131").Value = 10416.5
$ws.Range("K131").Value = 7497.500100000001
$ws.Range("L131").Value = 31249.5
$ws.Range("M131").Value = -2457.500100000001
$ws.Range("N131").Value = -41329.5

$ws.Range("H132").Value = 54245.95
$ws.Range("I132").Value = 30951.383
$ws.Range("K132").Value = 92854.149
$ws.Range("M132").Value = -90324.149

$ws.Range("H137").Value = 1454818.2
$ws.Range("I137").Value = 21243.768
$ws.Range("J137").Value = 7598708.5
$ws.Range("K137").Value = 63731.304
$ws.Range("L137").Value = 22796125.5
$ws.Range("M137").Value = -61181.304
$ws.Range("N137").Value = -22801225.5

$ws.Range("H138").Value = 3997.1685
$ws.Range("I138").Value = 2234.2307
$ws.Range("J138").Value = 4298.7236
$ws.Range("K138").Value = 6702.6921
$ws.Range("L138").Value = 12896.1708
$ws.Range("M138").Value = -1562.6921
$ws.Range("N138").Value = -23176.1708

$ws.Range("H139").Value = 119130
$ws.Range("J139").Value = 118695
$ws.Range("L139").Value = 118695
$ws.Range("N139").Value = -128975

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 6498.25
$ws.Range("I28").Value = 6498.25
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 6498.25
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -6306.25
$ws.Range("N28").ClearContents()

$ws.Range("H63").Value = 200
$ws.Range("I63").Value = 200
$ws.Range("K63").Value = 200
$ws.Range("M63").Value = 486

$ws.Range("H66").Value = 200
$ws.Range("I66").Value = 200
$ws.Range("K66").Value = 1000
$ws.Range("M66").Value = 2432

$ws.Range("H97").Value = 1285
$ws.Range("I97").Value = 492
$ws.Range("J97").Value = 5250
$ws.Range("K97").Value = 492
$ws.Range("L97").Value = 5250
$ws.Range("M97").Value = 4
$ws.Range("N97").Value = -6242

$ws.Range("H99").Value = 6498.25
$ws.Range("I99").Value = 6498.25
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 6498.25
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3503.25
$ws.Range("N99").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2437.923
$ws.Range("J94").Value = 2097.5
$ws.Range("L94").Value = 2097.5
$ws.Range("N94").Value = -2999.5

$ws.Range("H130").Value = 85000
$ws.Range("J130").Value = 85000
$ws.Range("L130").Value = 85000
$ws.Range("N130").Value = -95040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1274.2
$ws.Range("I7").Value = 64.25
$ws.Range("K7").Value = 64.25
$ws.Range("M7").Value = 48.75

$ws.Range("H58").Value = 2682.3157
$ws.Range("I58").Value = 2331
$ws.Range("K58").Value = 2331
$ws.Range("M58").Value = -2128

$ws.Range("H105").Value = 1671.1111
$ws.Range("I105").Value = 1671.1111
$ws.Range("K105").Value = 1671.1111
$ws.Range("M105").Value = 75.88889999999992

$ws.Range("H107").Value = 748.2727
$ws.Range("I107").Value = 820.5263
$ws.Range("J107").Value = 290.66666
$ws.Range("K107").Value = 820.5263
$ws.Range("L107").Value = 290.66666
$ws.Range("M107").Value = 1099.4737
$ws.Range("N107").Value = -4130.66666

$ws.Range("H136").Value = 2682.3157
$ws.Range("I136").Value = 2331
$ws.Range("K136").Value = 6993
$ws.Range("M136").Value = -4443

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 156.66667
$ws.Range("I23").Value = 224
$ws.Range("J23").Value = 72.5
$ws.Range("K23").Value = 672
$ws.Range("L23").Value = 217.5
$ws.Range("M23").Value = -437
$ws.Range("N23").Value = -687.5

$ws.Range("H68").Value = 1854508.9
$ws.Range("J68").Value = 1925663.8
$ws.Range("L68").Value = 5776991.4
$ws.Range("N68").Value = -5778613.4

$ws.Range("H71").Value = 1854508.9
$ws.Range("J71").Value = 1925663.8
$ws.Range("L71").Value = 17330974.2
$ws.Range("N71").Value = -17339086.2

$ws.Range("H86").Value = 215.11765
$ws.Range("I86").Value = 51.384617
$ws.Range("K86").Value = 154.153851
$ws.Range("M86").Value = 1031.846149

$ws.Range("H89").Value = 215.11765
$ws.Range("I89").Value = 51.384617
$ws.Range("K89").Value = 462.461553
$ws.Range("M89").Value = 5465.538447

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1225.25
$ws.Range("I9").Value = 1203
$ws.Range("J9").Value = 1247.5
$ws.Range("K9").Value = 1203
$ws.Range("L9").Value = 1247.5
$ws.Range("M9").Value = -1033
$ws.Range("N9").Value = -1587.5

$ws.Range("H80").Value = 2673.4211
$ws.Range("I80").Value = 2369.077
$ws.Range("J80").Value = 3332.8333
$ws.Range("K80").Value = 2369.077
$ws.Range("L80").Value = 3332.8333
$ws.Range("M80").Value = -1371.077
$ws.Range("N80").Value = -5328.8333

$ws.Range("H83").Value = 2673.4211
$ws.Range("I83").Value = 2369.077
$ws.Range("J83").Value = 3332.8333
$ws.Range("K83").Value = 11845.385
$ws.Range("L83").Value = 16664.1665
$ws.Range("M83").Value = -6853.385000000002
$ws.Range("N83").Value = -26648.1665

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H122").Value = 200
$ws.Range("I122").Value = 200
$ws.Range("K122").Value = 600
$ws.Range("M122").Value = 1850

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11578.571
$ws.Range("I7").Value = 14874.75
$ws.Range("J7").Value = 7183.6665
$ws.Range("K7").Value = 14874.75
$ws.Range("L7").Value = 7183.6665
$ws.Range("M7").Value = -14762.75
$ws.Range("N7").Value = -7407.6665

$ws.Range("H22").Value = 977.2727
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 1450
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 1450
$ws.Range("M22").Value = -505
$ws.Range("N22").Value = -2040

$ws.Range("H27").Value = 977.2727
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 1450
$ws.Range("K27").Value = 800
$ws.Range("L27").Value = 1450
$ws.Range("M27").Value = -693
$ws.Range("N27").Value = -1664

$ws.Range("H46").Value = 2444.7144
$ws.Range("I46").Value = 1461.9
$ws.Range("J46").Value = 3338.182
$ws.Range("K46").Value = 1461.9
$ws.Range("L46").Value = 3338.182
$ws.Range("M46").Value = -1273.9
$ws.Range("N46").Value = -3714.182

$ws.Range("H93").Value = 1016579.25
$ws.Range("I93").Value = 1856895.5
$ws.Range("K93").Value = 1856895.5
$ws.Range("M93").Value = -1855647.5

$ws.Range("H126").Value = 11578.571
$ws.Range("I126").Value = 14874.75
$ws.Range("J126").Value = 7183.6665
$ws.Range("K126").Value = 44624.25
$ws.Range("L126").Value = 21550.9995
$ws.Range("M126").Value = -42154.25
$ws.Range("N126").Value = -26490.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 19583.277
$ws.Range("I62").Value = 19705.854
$ws.Range("J62").Value = 17499.5
$ws.Range("K62").Value = 19705.854
$ws.Range("L62").Value = 17499.5
$ws.Range("M62").Value = -19081.854
$ws.Range("N62").Value = -18747.5

$ws.Range("H65").Value = 19583.277
$ws.Range("I65").Value = 19705.854
$ws.Range("J65").Value = 17499.5
$ws.Range("K65").Value = 98529.26999999999
$ws.Range("L65").Value = 87497.5
$ws.Range("M65").Value = -95409.26999999999
$ws.Range("N65").Value = -93737.5

$ws.Range("H92").Value = 40549.5
$ws.Range("J92").Value = 40549.5
$ws.Range("L92").Value = 40549.5
$ws.Range("N92").Value = -45541.5

$ws.Range("H107").Value = 1328.5
$ws.Range("I107").Value = 917.9474
$ws.Range("J107").Value = 2442.8572
$ws.Range("K107").Value = 2753.8422
$ws.Range("L107").Value = 7328.571599999999
$ws.Range("M107").Value = -833.8422
$ws.Range("N107").Value = -11168.5716

$ws.Range("H122").Value = 5122.3076
$ws.Range("I122").Value = 3167.75
$ws.Range("J122").Value = 8249.6
$ws.Range("K122").Value = 9503.25
$ws.Range("L122").Value = 24748.8
$ws.Range("M122").Value = -7053.25
$ws.Range("N122").Value = -29648.8

$ws.Range("H132").Value = 2032.2927
$ws.Range("I132").Value = 2230.7568
$ws.Range("K132").Value = 6692.2704
$ws.Range("M132").Value = -4162.2704

$ws.Range("H138").Value = 84899.5
$ws.Range("J138").Value = 79800
$ws.Range("L138").Value = 79800
$ws.Range("N138").Value = -90080
